$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the "percent_moisture" column (O2:O37) to the placeholder value
# used while the real moisture readings are recalculated.
$ws.Range("O2:O37").Value = 0.00001

# The header cell (O1) now carries the plain data-column font instead of
# the bordered/shaded header style.
$ws.Range("O1").Font.Name = "Calibri"

# The reset data cells (O2:O37) pick up the plain "Aptos Narrow" font
# (no border/fill) instead of the old bordered header-style font.
$ws.Range("O2:O37").Font.Name = "Aptos Narrow"

# Update the active selection to reflect the column being worked on.
$ws.Range("O1:O37").Select
